# Add a new worksheet "Resilience1" as the last sheet in the workbook and
# populate it with the resilience data / thesis table, matching the
# "resilience data and thesis updates" commit.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet after the current last sheet, so it becomes the
#     final (5th) tab and the active one -----------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Resilience1"

# --- Headers ---------------------------------------------------------------
$ws.Range("A1").Value = '"optimal" resilience for 1 node and 3 edges'
$ws.Range("A2").Value = "Time Steps"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3

# --- Data rows ---------------------------------------------------------------
$ws.Range("B3").Value = 231.69999999999899
$ws.Range("C3").Value = 231.69999999999899
$ws.Range("D3").Value = 261.69999999999902

$ws.Range("B4").Value = 104.7
$ws.Range("C4").Value = 93
$ws.Range("D4").Value = 134.69999999999999

$ws.Range("B5").Value = 35.4
$ws.Range("C5").Value = 75.099999999999994
$ws.Range("D5").Value = 93.7

$ws.Range("B6").Value = 35.4
$ws.Range("C6").Value = 40.099999999999902
$ws.Range("D6").Value = 63.7

$ws.Range("B7").Value = 35.4
$ws.Range("C7").Value = 24.2
$ws.Range("D7").Value = 38.199999999999903

$ws.Range("B8").Value = 35.4
$ws.Range("C8").Value = 24.2
$ws.Range("D8").Value = 25.5

# --- Column A width (matches the "bestFit" look of the original sheet) -----
$ws.Columns.Item(1).ColumnWidth = 18.1666666666667

# --- Match the selection left active on the sheet in the source file -------
$ws.Range("G24").Select()
